$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Decrement column E (剩余) by 1 for every data row (2-99), except row 36
# which is left unchanged, matching the supplied diff.
for ($r = 2; $r -le 99; $r++) {
    if ($r -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($r, 5)
    $current = $cell.Value()
    $cell.Value = $current - 1
}
